$d = $word.ActiveDocument

# The "Dheeraj Chand" header paragraph is immediately followed by the
# PROFESSIONAL SUMMARY heading; the short-resume contact line is missing.
# Insert it as a new, centered paragraph right after the name, using
# Find/Replace with a paragraph-mark (^p) so the new paragraph inherits
# the plain (non-bold/non-large) run formatting of a fresh paragraph
# rather than the name run's bold/28pt character formatting.
$range = $d.Content
[void]$range.Find.Execute(
    "Dheeraj Chand", $true, $false, $false, $false, $false, $true, 1, $false,
    "Dheeraj Chand^p202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX",
    2
)
